$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price / volume-label data to match the latest
# GitHub Actions scrape. D-column prices are stored as text, so a leading
# apostrophe is used to keep Excel from re-interpreting them as numbers
# (which would lose the exact decimal formatting).

$ws.Range("D2").Value = "'265.26"
$ws.Range("D3").Value = "'22.66"
$ws.Range("D5").Value = "'0.06146"
$ws.Range("D6").Value = "'3.579"
$ws.Range("D7").Value = "'6.664"
$ws.Range("D8").Value = "'1.339"
$ws.Range("D9").Value = "'0.8304"
$ws.Range("D10").Value = "'0.01354"
$ws.Range("D11").Value = "'0.1590"
$ws.Range("D12").Value = "'0.08184"
$ws.Range("D14").Value = "'0.03151"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09244"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.901"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("D17").Value = "'0.001722"
$ws.Range("D18").Value = "'0.04882"
$ws.Range("D19").Value = "'0.006213"
$ws.Range("D20").Value = "'0.005264"
$ws.Range("D21").Value = "'0.001090"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.767"
$ws.Range("D24").Value = "'2.322"
$ws.Range("D25").Value = "'0.3341"
$ws.Range("D26").Value = "'0.1238"
$ws.Range("D27").Value = "'0.0002681"
$ws.Range("D40").Value = "'0.04613"
$ws.Range("D41").Value = "'0.006968"
$ws.Range("D42").Value = "'0.1135"
$ws.Range("D43").Value = "'0.003603"
$ws.Range("D44").Value = "'0.01080"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"
$ws.Range("D45").Value = "'0.00006137"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.7895"
$ws.Range("D48").Value = "'0.1916"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "48CryptobidCoinCBC"
$ws.Range("D50").Value = "'0.01241"
